{"js": "// Replace the two-digit multiplication problems in the document body's\n// table cells with their updated values, matching the diff exactly.\n// Each \"old\" text value is unique in the document, so an exact,\n// case-sensitive, whole-text search/replace is unambiguous.\nconst replacements = [\n  [\"29\u00d725=\", \"33\u00d774=\"],\n  [\"74\u00d767=\", \"24\u00d742=\"],\n  [\"11\u00d771=\", \"62\u00d781=\"],\n  [\"56\u00d712=\", \"81\u00d799=\"],\n  [\"70\u00d721=\", \"63\u00d743=\"],\n  [\"54\u00d712=\", \"56\u00d725=\"],\n  [\"69\u00d725=\", \"87\u00d793=\"],\n  [\"19\u00d755=\", \"39\u00d772=\"],\n  [\"92\u00d770=\", \"78\u00d738=\"],\n  [\"43\u00d772=\", \"39\u00d740=\"],\n  [\"55\u00d747=\", \"95\u00d735=\"],\n  [\"40\u00d799=\", \"14\u00d725=\"],\n  [\"53\u00d732=\", \"53\u00d755=\"],\n  [\"92\u00d784=\", \"97\u00d736=\"],\n  [\"97\u00d726=\", \"57\u00d763=\"],\n  [\"83\u00d798=\", \"72\u00d776=\"],\n  [\"55\u00d757=\", \"51\u00d782=\"],\n  [\"57\u00d764=\", \"85\u00d744=\"],\n  [\"83\u00d742=\", \"97\u00d788=\"],\n  [\"45\u00d731=\", \"92\u00d796=\"],\n  [\"49\u00d725=\", \"17\u00d741=\"],\n  [\"21\u00d767=\", \"82\u00d725=\"],\n  [\"49\u00d793=\", \"62\u00d725=\"],\n  [\"22\u00d730=\", \"29\u00d766=\"],\n  [\"73\u00d775=\", \"22\u00d744=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the two-digit multiplication problems in the document's table\n# cells with their updated values, matching the diff exactly. Each \"old\"\n# text value is unique in the document, so an exact, whole-text\n# find/replace (wdReplaceAll) is unambiguous for every pair.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"29\u00d725=\", \"33\u00d774=\"),\n  @(\"74\u00d767=\", \"24\u00d742=\"),\n  @(\"11\u00d771=\", \"62\u00d781=\"),\n  @(\"56\u00d712=\", \"81\u00d799=\"),\n  @(\"70\u00d721=\", \"63\u00d743=\"),\n  @(\"54\u00d712=\", \"56\u00d725=\"),\n  @(\"69\u00d725=\", \"87\u00d793=\"),\n  @(\"19\u00d755=\", \"39\u00d772=\"),\n  @(\"92\u00d770=\", \"78\u00d738=\"),\n  @(\"43\u00d772=\", \"39\u00d740=\"),\n  @(\"55\u00d747=\", \"95\u00d735=\"),\n  @(\"40\u00d799=\", \"14\u00d725=\"),\n  @(\"53\u00d732=\", \"53\u00d755=\"),\n  @(\"92\u00d784=\", \"97\u00d736=\"),\n  @(\"97\u00d726=\", \"57\u00d763=\"),\n  @(\"83\u00d798=\", \"72\u00d776=\"),\n  @(\"55\u00d757=\", \"51\u00d782=\"),\n  @(\"57\u00d764=\", \"85\u00d744=\"),\n  @(\"83\u00d742=\", \"97\u00d788=\"),\n  @(\"45\u00d731=\", \"92\u00d796=\"),\n  @(\"49\u00d725=\", \"17\u00d741=\"),\n  @(\"21\u00d767=\", \"82\u00d725=\"),\n  @(\"49\u00d793=\", \"62\u00d725=\"),\n  @(\"22\u00d730=\", \"29\u00d766=\"),\n  @(\"73\u00d775=\", \"22\u00d744=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $newText\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
